# Edit script: apply Duplicate_Transactions diff
# 1) Fix B410: convert from text to a true numeric value (2065044242)
# 2) Append 22 new duplicate-transaction rows (411-432)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 410: Phone column should be numeric, not text ---
$ws.Cells.Item(410, 2).Value = 2065044242

# --- Append 22 new rows (411-432) recording duplicate-transaction attempts ---

# Row 411
$ws.Cells.Item(411, 1).Value = "David"
$ws.Cells.Item(411, 2).Value = 12814104622
$ws.Cells.Item(411, 3).Value = "PO Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(411, 4).Value = "GG"
$ws.Cells.Item(411, 5).Value = "Po Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(411, 11).Value = "WhatsApp message already sent for this book previously"
$ws.Cells.Item(411, 12).Value = "2025-09-19 12:50:13"
$ws.Cells.Item(411, 13).Value = "'2025-09-19"
$ws.Cells.Item(411, 14).Value = "Blocked"

# Row 412
$ws.Cells.Item(412, 1).Value = "Henry Chelegbor"
$ws.Cells.Item(412, 2).Value = 13024705411
$ws.Cells.Item(412, 3).Value = "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
$ws.Cells.Item(412, 4).Value = "GG"
$ws.Cells.Item(412, 5).Value = "6613 Guyer Street, Philadelphia, Pa, Pennsylvania"
$ws.Cells.Item(412, 11).Value = "WhatsApp message already sent for this book previously"
$ws.Cells.Item(412, 12).Value = "2025-09-19 12:50:14"
$ws.Cells.Item(412, 13).Value = "'2025-09-19"
$ws.Cells.Item(412, 14).Value = "Blocked"

# Row 413
$ws.Cells.Item(413, 1).Value = "Dennis Vanmeter"
$ws.Cells.Item(413, 2).Value = 13049196111
$ws.Cells.Item(413, 3).Value = "1909 Harper Rd, Beckley, WV 25801"
$ws.Cells.Item(413, 4).Value = "GG"
$ws.Cells.Item(413, 5).Value = "1909 Harper Rd, Beckley, Wv 25801"
$ws.Cells.Item(413, 11).Value = "WhatsApp message already sent for this book previously"
$ws.Cells.Item(413, 12).Value = "2025-09-19 12:50:16"
$ws.Cells.Item(413, 13).Value = "'2025-09-19"
$ws.Cells.Item(413, 14).Value = "Blocked"

# Row 414
$ws.Cells.Item(414, 1).Value = "Madhukar Verma"
$ws.Cells.Item(414, 2).Value = 2065044242
$ws.Cells.Item(414, 3).Value = "42729 Mayfair Park Ave Fremont Fremont 94538 California USA"
$ws.Cells.Item(414, 4).Value = "YBB"
$ws.Cells.Item(414, 5).Value = "English"
$ws.Cells.Item(414, 11).Value = "Same book already sent"
$ws.Cells.Item(414, 12).Value = "2025-09-19 12:53:09"
$ws.Cells.Item(414, 13).Value = "'2025-09-19"
$ws.Cells.Item(414, 14).Value = "Blocked"

# Row 415
$ws.Cells.Item(415, 1).Value = "David"
$ws.Cells.Item(415, 2).Value = 12814104622
$ws.Cells.Item(415, 3).Value = "PO Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(415, 4).Value = "GG"
$ws.Cells.Item(415, 5).Value = "Po Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(415, 11).Value = "Same book already sent"
$ws.Cells.Item(415, 12).Value = "2025-09-19 12:53:12"
$ws.Cells.Item(415, 13).Value = "'2025-09-19"
$ws.Cells.Item(415, 14).Value = "Blocked"

# Row 416
$ws.Cells.Item(416, 1).Value = "Henry Chelegbor"
$ws.Cells.Item(416, 2).Value = 13024705411
$ws.Cells.Item(416, 3).Value = "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
$ws.Cells.Item(416, 4).Value = "GG"
$ws.Cells.Item(416, 5).Value = "6613 Guyer Street, Philadelphia, Pa, Pennsylvania"
$ws.Cells.Item(416, 11).Value = "Same book already sent"
$ws.Cells.Item(416, 12).Value = "2025-09-19 12:53:14"
$ws.Cells.Item(416, 13).Value = "'2025-09-19"
$ws.Cells.Item(416, 14).Value = "Blocked"

# Row 417
$ws.Cells.Item(417, 1).Value = "Dennis Vanmeter"
$ws.Cells.Item(417, 2).Value = 13049196111
$ws.Cells.Item(417, 3).Value = "1909 Harper Rd, Beckley, WV 25801"
$ws.Cells.Item(417, 4).Value = "GG"
$ws.Cells.Item(417, 5).Value = "1909 Harper Rd, Beckley, Wv 25801"
$ws.Cells.Item(417, 11).Value = "Same book already sent"
$ws.Cells.Item(417, 12).Value = "2025-09-19 12:53:17"
$ws.Cells.Item(417, 13).Value = "'2025-09-19"
$ws.Cells.Item(417, 14).Value = "Blocked"

# Row 418
$ws.Cells.Item(418, 1).Value = "David"
$ws.Cells.Item(418, 2).Value = 12814104622
$ws.Cells.Item(418, 3).Value = "PO Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(418, 4).Value = "GG"
$ws.Cells.Item(418, 5).Value = "Po Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(418, 11).Value = "Same book already sent"
$ws.Cells.Item(418, 12).Value = "2025-09-19 12:53:49"
$ws.Cells.Item(418, 13).Value = "'2025-09-19"
$ws.Cells.Item(418, 14).Value = "Blocked"

# Row 419
$ws.Cells.Item(419, 1).Value = "Henry Chelegbor"
$ws.Cells.Item(419, 2).Value = 13024705411
$ws.Cells.Item(419, 3).Value = "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
$ws.Cells.Item(419, 4).Value = "GG"
$ws.Cells.Item(419, 5).Value = "6613 Guyer Street, Philadelphia, Pa, Pennsylvania"
$ws.Cells.Item(419, 11).Value = "Same book already sent"
$ws.Cells.Item(419, 12).Value = "2025-09-19 12:53:51"
$ws.Cells.Item(419, 13).Value = "'2025-09-19"
$ws.Cells.Item(419, 14).Value = "Blocked"

# Row 420
$ws.Cells.Item(420, 1).Value = "Dennis Vanmeter"
$ws.Cells.Item(420, 2).Value = 13049196111
$ws.Cells.Item(420, 3).Value = "1909 Harper Rd, Beckley, WV 25801"
$ws.Cells.Item(420, 4).Value = "GG"
$ws.Cells.Item(420, 5).Value = "1909 Harper Rd, Beckley, Wv 25801"
$ws.Cells.Item(420, 11).Value = "Same book already sent"
$ws.Cells.Item(420, 12).Value = "2025-09-19 12:53:53"
$ws.Cells.Item(420, 13).Value = "'2025-09-19"
$ws.Cells.Item(420, 14).Value = "Blocked"

# Row 421
$ws.Cells.Item(421, 1).Value = "David"
$ws.Cells.Item(421, 2).Value = 12814104622
$ws.Cells.Item(421, 3).Value = "PO Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(421, 4).Value = "GG"
$ws.Cells.Item(421, 5).Value = "Po Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(421, 11).Value = "Same book already sent"
$ws.Cells.Item(421, 12).Value = "2025-09-19 12:59:03"
$ws.Cells.Item(421, 13).Value = "'2025-09-19"
$ws.Cells.Item(421, 14).Value = "Blocked"

# Row 422
$ws.Cells.Item(422, 1).Value = "Henry Chelegbor"
$ws.Cells.Item(422, 2).Value = 13024705411
$ws.Cells.Item(422, 3).Value = "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
$ws.Cells.Item(422, 4).Value = "GG"
$ws.Cells.Item(422, 5).Value = "6613 Guyer Street, Philadelphia, Pa, Pennsylvania"
$ws.Cells.Item(422, 11).Value = "Same book already sent"
$ws.Cells.Item(422, 12).Value = "2025-09-19 12:59:05"
$ws.Cells.Item(422, 13).Value = "'2025-09-19"
$ws.Cells.Item(422, 14).Value = "Blocked"

# Row 423
$ws.Cells.Item(423, 1).Value = "Dennis Vanmeter"
$ws.Cells.Item(423, 2).Value = 13049196111
$ws.Cells.Item(423, 3).Value = "1909 Harper Rd, Beckley, WV 25801"
$ws.Cells.Item(423, 4).Value = "GG"
$ws.Cells.Item(423, 5).Value = "1909 Harper Rd, Beckley, Wv 25801"
$ws.Cells.Item(423, 11).Value = "Same book already sent"
$ws.Cells.Item(423, 12).Value = "2025-09-19 12:59:07"
$ws.Cells.Item(423, 13).Value = "'2025-09-19"
$ws.Cells.Item(423, 14).Value = "Blocked"

# Row 424
$ws.Cells.Item(424, 1).Value = "David"
$ws.Cells.Item(424, 2).Value = 12814104622
$ws.Cells.Item(424, 3).Value = "PO Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(424, 4).Value = "GG"
$ws.Cells.Item(424, 5).Value = "Po Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(424, 11).Value = "Same book already sent"
$ws.Cells.Item(424, 12).Value = "2025-09-19 13:00:34"
$ws.Cells.Item(424, 13).Value = "'2025-09-19"
$ws.Cells.Item(424, 14).Value = "Blocked"

# Row 425
$ws.Cells.Item(425, 1).Value = "Henry Chelegbor"
$ws.Cells.Item(425, 2).Value = 13024705411
$ws.Cells.Item(425, 3).Value = "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
$ws.Cells.Item(425, 4).Value = "GG"
$ws.Cells.Item(425, 5).Value = "6613 Guyer Street, Philadelphia, Pa, Pennsylvania"
$ws.Cells.Item(425, 11).Value = "Same book already sent"
$ws.Cells.Item(425, 12).Value = "2025-09-19 13:00:36"
$ws.Cells.Item(425, 13).Value = "'2025-09-19"
$ws.Cells.Item(425, 14).Value = "Blocked"

# Row 426
$ws.Cells.Item(426, 1).Value = "Dennis Vanmeter"
$ws.Cells.Item(426, 2).Value = 13049196111
$ws.Cells.Item(426, 3).Value = "1909 Harper Rd, Beckley, WV 25801"
$ws.Cells.Item(426, 4).Value = "GG"
$ws.Cells.Item(426, 5).Value = "1909 Harper Rd, Beckley, Wv 25801"
$ws.Cells.Item(426, 11).Value = "Same book already sent"
$ws.Cells.Item(426, 12).Value = "2025-09-19 13:00:39"
$ws.Cells.Item(426, 13).Value = "'2025-09-19"
$ws.Cells.Item(426, 14).Value = "Blocked"

# Row 427
$ws.Cells.Item(427, 1).Value = "David"
$ws.Cells.Item(427, 2).Value = 12814104622
$ws.Cells.Item(427, 3).Value = "PO Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(427, 4).Value = "GG"
$ws.Cells.Item(427, 5).Value = "Po Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(427, 11).Value = "WhatsApp message already sent for this book previously"
$ws.Cells.Item(427, 12).Value = "2025-09-19 13:04:39"
$ws.Cells.Item(427, 13).Value = "'2025-09-19"
$ws.Cells.Item(427, 14).Value = "Blocked"

# Row 428
$ws.Cells.Item(428, 1).Value = "Henry Chelegbor"
$ws.Cells.Item(428, 2).Value = 13024705411
$ws.Cells.Item(428, 3).Value = "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
$ws.Cells.Item(428, 4).Value = "GG"
$ws.Cells.Item(428, 5).Value = "6613 Guyer Street, Philadelphia, Pa, Pennsylvania"
$ws.Cells.Item(428, 11).Value = "WhatsApp message already sent for this book previously"
$ws.Cells.Item(428, 12).Value = "2025-09-19 13:04:41"
$ws.Cells.Item(428, 13).Value = "'2025-09-19"
$ws.Cells.Item(428, 14).Value = "Blocked"

# Row 429
$ws.Cells.Item(429, 1).Value = "Dennis Vanmeter"
$ws.Cells.Item(429, 2).Value = 13049196111
$ws.Cells.Item(429, 3).Value = "1909 Harper Rd, Beckley, WV 25801"
$ws.Cells.Item(429, 4).Value = "GG"
$ws.Cells.Item(429, 5).Value = "1909 Harper Rd, Beckley, Wv 25801"
$ws.Cells.Item(429, 11).Value = "WhatsApp message already sent for this book previously"
$ws.Cells.Item(429, 12).Value = "2025-09-19 13:04:43"
$ws.Cells.Item(429, 13).Value = "'2025-09-19"
$ws.Cells.Item(429, 14).Value = "Blocked"

# Row 430
$ws.Cells.Item(430, 1).Value = "David"
$ws.Cells.Item(430, 2).Value = 12814104622
$ws.Cells.Item(430, 3).Value = "PO Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(430, 4).Value = "GG"
$ws.Cells.Item(430, 5).Value = "Po Box 87301, Park Place, Houston, Texas"
$ws.Cells.Item(430, 11).Value = "Same book already sent"
$ws.Cells.Item(430, 12).Value = "2025-09-19 13:05:43"
$ws.Cells.Item(430, 13).Value = "'2025-09-19"
$ws.Cells.Item(430, 14).Value = "Blocked"

# Row 431
$ws.Cells.Item(431, 1).Value = "Henry Chelegbor"
$ws.Cells.Item(431, 2).Value = 13024705411
$ws.Cells.Item(431, 3).Value = "6613 Guyer Street, Philadelphia, PA, Pennsylvania"
$ws.Cells.Item(431, 4).Value = "GG"
$ws.Cells.Item(431, 5).Value = "6613 Guyer Street, Philadelphia, Pa, Pennsylvania"
$ws.Cells.Item(431, 11).Value = "Same book already sent"
$ws.Cells.Item(431, 12).Value = "2025-09-19 13:05:45"
$ws.Cells.Item(431, 13).Value = "'2025-09-19"
$ws.Cells.Item(431, 14).Value = "Blocked"

# Row 432
$ws.Cells.Item(432, 1).Value = "Dennis Vanmeter"
$ws.Cells.Item(432, 2).Value = "'13049196111"
$ws.Cells.Item(432, 3).Value = "1909 Harper Rd, Beckley, WV 25801"
$ws.Cells.Item(432, 4).Value = "GG"
$ws.Cells.Item(432, 5).Value = "1909 Harper Rd, Beckley, Wv 25801"
$ws.Cells.Item(432, 11).Value = "Same book already sent"
$ws.Cells.Item(432, 12).Value = "2025-09-19 13:05:47"
$ws.Cells.Item(432, 13).Value = "'2025-09-19"
$ws.Cells.Item(432, 14).Value = "Blocked"
